$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.336.35'
$ws.Range("E2").Value = '  -0.54%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.497.11'
$ws.Range("E3").Value = '  -0.54%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.04'
$ws.Range("E5").Value = '  -0.10%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.99'
$ws.Range("E6").Value = '  +1.73%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.496.27'
$ws.Range("E7").Value = '  -0.60%  '

$ws.Range("E8").Value = '  +0.05%  '

$ws.Range("E9").Value = '  -0.85%  '

$ws.Range("E10").Value = '  -0.27%  '

$ws.Range("E11").Value = '  -1.87%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.375'
$ws.Range("E12").Value = '  -3.44%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.095.68'

$ws.Range("E14").Value = '  +0.12%  '

$ws.Range("E15").Value = '  +1.26%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.501.32'
$ws.Range("E16").Value = '  -0.31%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.376.06'
$ws.Range("E17").Value = '  -0.38%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.03'
$ws.Range("E18").Value = '  -9.89%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.00'
$ws.Range("E19").Value = '  +0.16%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.64'
$ws.Range("E20").Value = '  -0.75%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.75'
$ws.Range("E21").Value = '  -3.22%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '384.64'
$ws.Range("E22").Value = '  -1.77%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.568'
$ws.Range("E23").Value = '  -1.94%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.639.48'
$ws.Range("E24").Value = '  -0.44%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.06'
$ws.Range("E25").Value = '  +0.29%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.12%  '

$ws.Range("E27").Value = '  +2.80%  '

$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.47'
$ws.Range("E28").Value = '  -0.26%  '

$ws.Range("B29").Value = 'Fetch.AI'
$ws.Range("C29").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.55'
$ws.Range("E29").Value = '  -2.39%  '

$ws.Range("E30").Value = '  +0.07%  '

$ws.Range("E31").Value = '  -1.47%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.23'
$ws.Range("E32").Value = '  +0.04%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.518.33'
$ws.Range("E33").Value = '  +0.09%  '

$ws.Range("E34").Value = '  -0.03%  '

$ws.Range("E35").Value = '  +0.77%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '23.53'
$ws.Range("E36").Value = '  -2.13%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.28'
$ws.Range("E37").Value = '  +0.05%  '

$ws.Range("E38").Value = '  -3.56%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.82'
$ws.Range("E39").Value = '  -2.42%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '162.71'
$ws.Range("E40").Value = '  -4.84%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0780'
$ws.Range("E41").Value = '  -3.31%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.806'
$ws.Range("E42").Value = '  -1.10%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '25.81'
$ws.Range("E43").Value = '  -3.21%  '

$ws.Range("E44").Value = '  +0.10%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.89'
$ws.Range("E45").Value = '  -0.61%  '

$ws.Range("E46").Value = '  -0.76%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.40'
$ws.Range("E47").Value = '  +0.01%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.65'
$ws.Range("E48").Value = '  -0.67%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.477.62'
$ws.Range("E49").Value = '  +1.14%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.76'
$ws.Range("E50").Value = '  -1.99%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.908'
$ws.Range("E51").Value = '  +0.57%  '
